# Auto-generated Excel COM-interop script applying the cryptos.xlsx crypto-price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin) and C (Link) are reordered/renamed via plain text assignment --
# these values never look like numbers, so Excel keeps them as text automatically.
# Column D (Price) and E (Volume) values are always stored as *text* in this sheet
# (several already contain two dots, e.g. "76.413.39", which cannot be numbers).
# Assigning a plain numeric-looking string via .Value would make Excel coerce it to
# a real number, so for every D/E cell we briefly force a text format, assign the
# value, then clear the formatting again so the cell keeps its original (default)
# style while the stored value remains text, matching the source workbook exactly.

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '76.420.83'
Set-TextValue $ws.Range("E2") '  +0.59%  '
Set-TextValue $ws.Range("D3") '3.030.34'
Set-TextValue $ws.Range("E3") '  +4.19%  '
Set-TextValue $ws.Range("E4") '  +0.01%  '
Set-TextValue $ws.Range("D5") '199.88'
Set-TextValue $ws.Range("E5") '  +0.69%  '
Set-TextValue $ws.Range("D6") '628.49'
Set-TextValue $ws.Range("E6") '  +4.82%  '
Set-TextValue $ws.Range("E8") '  +0.40%  '
Set-TextValue $ws.Range("E9") '  +1.95%  '
Set-TextValue $ws.Range("D10") '3.031.51'
Set-TextValue $ws.Range("E10") '  +4.32%  '
Set-TextValue $ws.Range("D11") '0.436'
Set-TextValue $ws.Range("E11") '  +0.84%  '
Set-TextValue $ws.Range("E12") '  -0.45%  '
Set-TextValue $ws.Range("E13") '  +4.30%  '
Set-TextValue $ws.Range("D14") '3.586.32'
Set-TextValue $ws.Range("E14") '  +4.17%  '
Set-TextValue $ws.Range("D15") '29.15'
Set-TextValue $ws.Range("E15") '  +6.46%  '
Set-TextValue $ws.Range("D16") '76.350.12'
Set-TextValue $ws.Range("E16") '  +0.73%  '
Set-TextValue $ws.Range("E17") '  -1.85%  '
Set-TextValue $ws.Range("D18") '3.034.81'
Set-TextValue $ws.Range("E18") '  +4.59%  '
Set-TextValue $ws.Range("D19") '13.43'
Set-TextValue $ws.Range("E19") '  +3.49%  '
Set-TextValue $ws.Range("D20") '9.03'
Set-TextValue $ws.Range("E20") '  +3.19%  '
Set-TextValue $ws.Range("D21") '372.24'
Set-TextValue $ws.Range("E21") '  +0.51%  '
Set-TextValue $ws.Range("D22") '4.36'
Set-TextValue $ws.Range("E22") '  +1.65%  '
Set-TextValue $ws.Range("E23") '  -1.77%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D24") '73.26'
Set-TextValue $ws.Range("E24") '  +3.00%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range("D25") '3.206.47'
Set-TextValue $ws.Range("E25") '  +4.77%  '
Set-TextValue $ws.Range("E26") '  -0.05%  '
Set-TextValue $ws.Range("E27") '  +4.04%  '
Set-TextValue $ws.Range("D28") '9.91'
Set-TextValue $ws.Range("E28") '  +2.73%  '
Set-TextValue $ws.Range("E29") '  -0.09%  '
Set-TextValue $ws.Range("D30") '0.997'
Set-TextValue $ws.Range("E30") '  +0.27%  '
Set-TextValue $ws.Range("D31") '8.29'
Set-TextValue $ws.Range("E31") '  +7.65%  '
Set-TextValue $ws.Range("D32") '1.41'
Set-TextValue $ws.Range("E32") '  -0.18%  '
Set-TextValue $ws.Range("D33") '508.83'
Set-TextValue $ws.Range("E33") '  +1.29%  '
Set-TextValue $ws.Range("D34") '1.94'
Set-TextValue $ws.Range("E34") '  +7.10%  '
Set-TextValue $ws.Range("D35") '0.999'
Set-TextValue $ws.Range("E35") '  -0.02%  '
Set-TextValue $ws.Range("D36") '20.71'
Set-TextValue $ws.Range("E36") '  +2.44%  '
Set-TextValue $ws.Range("D37") '164.03'
Set-TextValue $ws.Range("E37") '  -0.64%  '
Set-TextValue $ws.Range("D38") '193.62'
Set-TextValue $ws.Range("E38") '  +7.34%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D39") '20.01'
Set-TextValue $ws.Range("E39") '  +1.91%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D40") '0.383'
Set-TextValue $ws.Range("E40") '  +10.94%  '
Set-TextValue $ws.Range("E41") '  -0.49%  '
Set-TextValue $ws.Range("E42") '  -1.18%  '
Set-TextValue $ws.Range("E43") '  +0.35%  '
Set-TextValue $ws.Range("D44") '5.08'
Set-TextValue $ws.Range("E44") '  +1.87%  '
Set-TextValue $ws.Range("D45") '42.57'
Set-TextValue $ws.Range("E45") '  +6.27%  '
Set-TextValue $ws.Range("E46") '  +5.32%  '
Set-TextValue $ws.Range("E47") '  +0.33%  '
Set-TextValue $ws.Range("D48") '0.718'
Set-TextValue $ws.Range("E48") '  +9.14%  '
Set-TextValue $ws.Range("D49") '0.602'
Set-TextValue $ws.Range("E49") '  +5.46%  '
Set-TextValue $ws.Range("D50") '2.34'
Set-TextValue $ws.Range("E50") '  +0.48%  '
Set-TextValue $ws.Range("D51") '3.87'
Set-TextValue $ws.Range("E51") '  +3.94%  '
